$d = $word.ActiveDocument

# --- Edit 1: Professional summary paragraph -------------------------------
$d.Content.Find.Execute(
    "Discovered systematic demographic coding errors affecting all Black and Asian-American voters, developed geospatial ML",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Discovered systematic demographic coding errors affecting 50M voters, developed geospatial ML",
    2) | Out-Null

# --- Edit 2: Bullet point - split run so "50M" is bold + colored ----------
# First normalize the wording in plain text (keeps the run's existing
# formatting, which is unformatted for this whole sentence).
$d.Content.Find.Execute(
    "Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Discovered systematic race coding errors affecting 50M voters, developed geospatial machine learning",
    2) | Out-Null

# Now re-find just the "50M" token inside that bullet and give it the same
# bold + accent-color formatting used by the other highlighted figures
# (e.g. 23%, 64%) in this document.
$r2 = $d.Content
$r2.Find.Execute("50M voters, developed geospatial machine learning") | Out-Null
$r2.Collapse(1)
$r2.MoveEnd(1, 3) | Out-Null
$r2.Font.Bold = 1
$r2.Font.Color = 5258796

# --- Edit 3: Project impact statement --------------------------------------
$d.Content.Find.Execute(
    "Impact: Corrected demographic data affecting all Black and Asian-American voters, improved electoral prediction accuracy by 22%",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Impact: Corrected demographic data affecting 50M voters nationwide, improved electoral prediction accuracy by 22%",
    2) | Out-Null
